$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("unitario")

$ws.Rows.Item(18).Insert()

$ws.Cells.Item(18, 1).Value = "216000"
$ws.Cells.Item(18, 2).Value = "TROCA DE HIDROMETRO PREVENTIVA AGENDADA"
$ws.Cells.Item(18, 3).Value = "Hidrometro"

$newRow = $ws.Range("A18:C18")
$newRow.Font.Bold = $true
$newRow.Font.Color = 0
$newRow.Interior.Color = 5296274
$newRow.Borders.LineStyle = 1
$newRow.Borders.Color = 13027014
$ws.Range("A18").NumberFormat = "@"

$ws.Activate()
$ws.Range("A10").Select()
$ws.Range("D15").Select()
